$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.90"
$ws.Range("E2").Value = "'-0.20%"
$ws.Range("D3").Value = "'40.67"
$ws.Range("E3").Value = "'3.84%"
$ws.Range("D4").Value = "'5.109"
$ws.Range("E4").Value = "'2.14%"
$ws.Range("D5").Value = "'0.07600"
$ws.Range("E5").Value = "'-1.63%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.271"
$ws.Range("E6").Value = "'-0.13%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.619"
$ws.Range("E7").Value = "'1.66%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.454"
$ws.Range("E8").Value = "'-4.10%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9054"
$ws.Range("E9").Value = "'-1.29%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1011"
$ws.Range("E10").Value = "'0.77%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1748"
$ws.Range("E11").Value = "'0.84%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09155"
$ws.Range("E12").Value = "'1.42%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04273"
$ws.Range("E13").Value = "'-4.71%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1056"
$ws.Range("E14").Value = "'-0.36%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001246"
$ws.Range("E15").Value = "'-0.79%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005843"
$ws.Range("E16").Value = "'3.41%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.350"
$ws.Range("E17").Value = "'-0.40%"
$ws.Range("E18").Value = "'-2.77%"
$ws.Range("D19").Value = "'6.577"
$ws.Range("E19").Value = "'-6.66%"
$ws.Range("D20").Value = "'0.1356"
$ws.Range("E20").Value = "'-0.48%"
$ws.Range("D21").Value = "'0.2727"
$ws.Range("E21").Value = "'-4.77%"
$ws.Range("D22").Value = "'0.04172"
$ws.Range("E22").Value = "'0.60%"
$ws.Range("D23").Value = "'0.001231"
$ws.Range("E23").Value = "'2.58%"
$ws.Range("D24").Value = "'0.004072"
$ws.Range("E24").Value = "'-0.18%"
$ws.Range("D25").Value = "'0.0001302"
$ws.Range("D26").Value = "'0.0003012"
$ws.Range("E26").Value = "'0.68%"
$ws.Range("D38").Value = "'0.02369"
$ws.Range("E38").Value = "'1.26%"
$ws.Range("D39").Value = "'0.05134"
$ws.Range("E39").Value = "'0.08%"
$ws.Range("D40").Value = "'0.007785"
$ws.Range("E40").Value = "'-2.02%"
$ws.Range("D41").Value = "'0.1295"
$ws.Range("E41").Value = "'-2.46%"
$ws.Range("D42").Value = "'0.007075"
$ws.Range("E42").Value = "'-3.43%"
$ws.Range("D43").Value = "'0.001860"
$ws.Range("E43").Value = "'-6.39%"
$ws.Range("D44").Value = "'0.008449"
$ws.Range("E44").Value = "'5.37%"
$ws.Range("D45").Value = "'0.3315"
$ws.Range("E45").Value = "'0.16%"
$ws.Range("D46").Value = "'0.00006363"
$ws.Range("E46").Value = "'-4.91%"
$ws.Range("E47").Value = "'-0.27%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.007209"
$ws.Range("E48").Value = "'112.21%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.004407"
$ws.Range("E49").Value = "'7.03%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.27%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.27%"
